# "adding more data and updating links" — append a new resource row
# (Covidtracking.com) below the existing "BLS Time Series" row on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D10 (URL) before C10 (label) so the new shared-string entries land
# in the same order as the source workbook: URL first, then label text.
$ws.Range("D10").Value = "https://covidtracking.com/data/national"
$ws.Range("C10").Value = "Covidtracker - aggregated case data"

# Leave the selection where the author left it after entering the data.
$ws.Range("C11").Select() | Out-Null
